$d = $word.ActiveDocument

function Set-ParaText($idx, $newText, $sizePt) {
    $p = $d.Paragraphs.Item($idx)
    $s = $p.Range.Start
    $e = $p.Range.End
    if ($e -gt $s) {
        $r = $d.Range($s, $e - 1)
    } else {
        $r = $d.Range($s, $e)
    }
    $r.Text = $newText
    $ne = $s + $newText.Length
    $r2 = $d.Range($s, $ne)
    $r2.Font.Name = "Times New Roman"
    $r2.Font.Color = 0
    if ($sizePt -ne $null) {
        $r2.Font.Size = $sizePt
    }
}

# --- Pass 1: normalize font across every existing paragraph (excluding paragraph marks) ---
foreach ($p in $d.Paragraphs) {
    $ps = $p.Range.Start
    $pe = $p.Range.End
    if ($pe -gt $ps) {
        $fr = $d.Range($ps, $pe - 1)
        $fr.Font.Name = "Times New Roman"
    }
}

# --- Pass 2: content edits ---

# Title
Set-ParaText 1 'Exploring the Wonders of Chemical Reactions: Unraveling the Secrets of Matter''s Transformations' 22

# Author byline
Set-ParaText 2 'Dr. Ella Harrison' 18

# Email line
Set-ParaText 3 'dr.ellaharrison@gmail.com' 16

# Main body paragraph
$bodyText = 'Chemistry, the science that unravels the intricate world of matter and its transformations, stands as a captivating field of study. It delves into the fundamental principles governing chemical reactions, unveiling the mysteries behind the interactions of atoms and molecules. Imagine embarking on a fascinating journey through the world of chemistry, exploring the interplay of elements, and witnessing the marvels of chemical reactions firsthand. As we delve into the realm of chemistry, we will uncover the secrets of matter''s transformations, unraveling the enigmatic tapestry of chemical reactions.' + [char]11 + [char]11 + 'Chemical reactions, the very essence of chemistry, are the processes through which atoms and molecules rearrange themselves to form new substances. These transformations are driven by the relentless pursuit of stability, as atoms seek to achieve their lowest energy configurations. From the spectacular combustion of fuels to the intricate processes of photosynthesis, chemical reactions orchestrate a symphony of changes that shape our world. By understanding these reactions, we gain insights into the fundamental mechanisms that govern the universe itself.' + [char]11 + [char]11 + 'The exploration of chemistry is not merely an academic pursuit; it is a gateway to unraveling the complexities of the natural world. From the delicate balance of ecosystems to the intricate workings of the human body, chemistry plays a crucial role in unraveling the mysteries that surround us. As we delve into the fascinating realm of chemical reactions, we not only gain knowledge but also develop critical thinking skills, problem-solving abilities, and an appreciation for the interconnectedness of all things.'
Set-ParaText 5 $bodyText 12

# Summary body paragraph (no explicit font size -> inherits default)
Set-ParaText 7 'Chemistry, the science that explores the transformations of matter through chemical reactions, stands as a captivating field of study. Chemical reactions, driven by the pursuit of stability, orchestrate a symphony of changes that shape our world. From combustion to photosynthesis, these reactions provide insights into the fundamental mechanisms that govern the universe. By unraveling the secrets of chemical reactions, we gain knowledge, develop critical thinking skills, and appreciate the interconnectedness of all things. Chemistry is not just an academic pursuit; it is a gateway to comprehending the complexities of the natural world and the intricacies of life itself.' $null

# --- Pass 3: append trailing empty paragraph before the section break ---
$d.Content.InsertParagraphAfter()

Write-Output "OK"